$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Update the date field result "10 June 2016" -> "17 June 2016"
# ------------------------------------------------------------------
$find = $d.Content.Find
$find.Execute("10 June 2016", $true, $false, $false, $false, $false, $true, 1, $false, "17 June 2016", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Remove the old "_GoBack" bookmark that sits right after the
#    "??" run (figure caption placeholder). Removing it renumbers
#    every following bookmark down by one, which is exactly what
#    turns the "_Toc453234639" bookmark id from 6 into 5.
# ------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 3. Rewrite the "technologies used" paragraph:
#      - "are going to be used" -> "are used"
#      - add "JQuery, " as its own run before "etc.)."
#      - re-insert a "_GoBack" bookmark right after "JQuery, "
#    This is done by replacing the whole paragraph's WordOpenXML so
#    the exact run/bookmark layout from the target document is
#    reproduced faithfully.
# ------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*technologies used*") {
        $target = $p.Range
        break
    }
}

$xmlFrag = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="7AD62DD4" w14:textId="79701284" w:rsidR="008300D6" w:rsidRPr="008A71C8" w:rsidRDefault="00B825A1" w:rsidP="008A71C8"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:bCs/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-GB"/></w:rPr><w:t>The technologies used are</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-GB"/></w:rPr><w:t>html</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> for the GUI, CSS for the style and Javascript for the interaction. Furthermore, cartographic dedicated libraries are used (Leaflet D3, Openlayers, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">JQuery, </w:t></w:r><w:bookmarkStart w:id="6" w:name="_GoBack"/><w:bookmarkEnd w:id="6"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-GB"/></w:rPr><w:t>etc.).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.InsertXML($xmlFrag)
